$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.127.05'
$ws.Range('E2').Value = '  -1.39%  '
$ws.Range('D3').Value = '1.656.76'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5181'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2626'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06260'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07707'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.422'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = '1.646.94'
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').Value = '1.884.26'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5407'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.45%  '
$ws.Range('D16').Value = '0.0₅8110'
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.68'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '26.170.00'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.616'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.027'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.96%  '
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1226'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.169'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.400'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05947'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.267'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.537'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.254'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.98%  '
$ws.Range('E34').Value = '  -5.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9634'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.426'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.772'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5674'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01591'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.960'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8543'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').Value = '1.010.04'
$ws.Range('E43').Value = '  -7.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.41'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').Value = '1.798.99'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('D46').Value = '0.0₈108'
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.51'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.89%  '
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.962'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05172'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('E51').Value = '  -0.87%  '
